# Walsin resistor rule modified
# Replace the contents of the RESISTOR_PREFIX table (A1:E6 -> A1:E10)
# with the updated set of Walsin size-code / power-rating rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("RESISTOR_PREFIX")

# Clear out the old table range in case the new range is smaller/larger
# than the old one (it's larger here, so this is mostly precautionary).
$ws.Range("A1:E10").ClearContents()

$data = @(
    @("Prefix", "Rating_Value", "Rating_Unit", "Vendor", "Priority"),
    @("WR25X", 1,      "W", "WALSIN", 1),
    @("WR20X", 0.75,   "W", "WALSIN", 1),
    @("WR18",  1,      "W", "WALSIN", 1),
    @("WR10X", 0.5,    "W", "WALSIN", 1),
    @("WR12X", 0.25,   "W", "WALSIN", 1),
    @("WR08X", 0.125,  "W", "WALSIN", 1),
    @("WR06X", 0.1,    "W", "WALSIN", 1),
    @("WR04X", 0.0625, "W", "WALSIN", 1),
    @("WR02X", 0.05,   "W", "WALSIN", 1)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
